$d = $word.ActiveDocument

# Paragraph 1: title/heading
$r1 = $d.Paragraphs(1).Range
$r1.MoveEnd(1, -1) | Out-Null
$r1.Text = '"Navigating the Linguistic Realm: Unleashing the Power of Data Science in NLP"'

# Paragraph 2: intro paragraph
$r2 = $d.Paragraphs(2).Range
$r2.MoveEnd(1, -1) | Out-Null
$r2.Text = 'Natural Language Processing (NLP) is at the forefront of the rapidly advancing field of data science. With the exponential growth of digital information, the ability to extract insights from text data has become crucial for businesses and organizations. NLP combines computational linguistics with artificial intelligence techniques to enable machines to understand and process human language. By analyzing and interpreting text data, NLP allows us to uncover patterns, sentiments, and extract meaningful information from vast amounts of unstructured text. This paper explores the fascinating world of NLP and its applications in various domains, highlighting the significant impact it has on data science.'

# Paragraph 3: body paragraph (multiple segments separated by manual line breaks)
$r3 = $d.Paragraphs(3).Range
$r3.MoveEnd(1, -1) | Out-Null
$r3.Text = 'Natural Language Processing (NLP) is a subfield of artificial intelligence and data science that focuses on the interaction between computers and human language. It combines principles from computer science, linguistics, and statistics to enable machines to understand, interpret, and analyze human language in a meaningful way. NLP plays a crucial role in various applications, from information retrieval and sentiment analysis to machine translation and speech recognition. Its potential is further enhanced when combined with data science techniques, opening up new possibilities for the analysis and processing of large volumes of textual data.' + [char]11 + [char]11 + 'One of the key challenges in NLP is to enable computers to understand human language, which is inherently complex and ambiguous. Data science techniques, such as machine learning and statistical modeling, have been instrumental in addressing this challenge. These techniques enable computers to learn patterns and relationships from large datasets, allowing them to recognize and interpret various aspects of language, such as the meaning of words and phrases, syntactic structure, and semantic relationships.' + [char]11 + [char]11 + 'Machine learning algorithms, such as neural networks and decision trees, have been used extensively in NLP to extract insights and patterns from text data. For example, in sentiment analysis, machine learning models can be trained on labeled datasets to classify text as positive, negative, or neutral. By analyzing large volumes of textual data, these models can provide valuable insights into people''s opinions and sentiments towards products, services, or events.' + [char]11 + [char]11 + 'Another important aspect of NLP is information extraction, which involves identifying and extracting structured information from unstructured text. Data science techniques, such as named entity recognition and relation extraction, play a crucial role in this process. For example, in a news article, data scientists can use NLP algorithms to identify and extract entities such as people, organizations, and locations, as well as relationships between them. This information can then be used for various purposes, such as building knowledge graphs or analyzing social networks.' + [char]11 + [char]11 + 'Furthermore, NLP and data science are closely intertwined in the field of text generation and language modeling. Language models, such as recurrent neural networks and transformers, have revolutionized the way computers generate human-like text. These models are trained on large datasets, learning to capture the statistical patterns and structures of language. By leveraging these models, data scientists can generate realistic and coherent text for a wide range of applications, from chatbots and virtual assistants to content generation and storytelling.' + [char]11 + [char]11 + 'In addition to text analysis and generation, NLP and data science techniques are also used in speech recognition and natural language understanding. Speech recognition algorithms convert spoken language into text, enabling machines to understand and process spoken commands or queries. This is particularly useful in applications such as voice assistants and speech-to-text transcription. Similarly, natural language understanding algorithms enable machines to interpret and respond to human queries in a meaningful way, allowing for more intuitive interactions between humans and computers.' + [char]11 + [char]11 + 'The field of NLP and data science is constantly evolving, with new techniques and methodologies being developed to tackle complex language-related challenges. As more and more data becomes available and computational power increases, the potential of NLP and data science in understanding and processing human language will continue to grow. This has significant implications in areas such as healthcare, customer service, social media analysis, and language translation.' + [char]11 + [char]11 + 'In conclusion, the combination of NLP and data science has revolutionized the way computers understand, interpret, and analyze human language. Data science techniques enable machines to learn patterns and relationships from large textual datasets, allowing for accurate language processing. From sentiment analysis to information extraction and text generation, NLP and data science have numerous applications that can benefit a wide range of industries. As technology continues to advance, the potential of NLP and data science in understanding and processing human language will continue to expand, opening up new possibilities for innovation and discovery.'

# Paragraph 4: conclusion paragraph
$r4 = $d.Paragraphs(4).Range
$r4.MoveEnd(1, -1) | Out-Null
$r4.Text = 'In conclusion, the field of Natural Language Processing (NLP) holds immense potential for data science applications. By combining the power of machine learning and linguistic analysis, NLP enables computers to understand, interpret, and generate human language. With the exponential growth of data and the increasing need to extract valuable insights from text, NLP is becoming integral to various industries, including healthcare, finance, marketing, and customer service. As advancements in NLP continue to push boundaries, we can expect to see innovative solutions that enhance automation, improve decision making, and revolutionize the way we interact with technology. Embracing the power of NLP is crucial in unlocking the full potential of data science and paving the way for a smarter, more connected future.'
